$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.206.41'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.790.30'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4529'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +20.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.136'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07457'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.189'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.221'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").Value = '1.792.36'
$ws.Range("E16").Value = '  +1.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001079'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06675'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '80.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.364'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").Value = '28.203.70'
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.385'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.372'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("D29").Value = '1.997.21'
$ws.Range("E29").Value = '  +1.75%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.260'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.071'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.855'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09388'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02362'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6607'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06223'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.155'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2147'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("E41").Value = '  +2.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.210'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.043'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.862'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6043'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.013'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07069'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("E51").Value = '  -1.57%  '
